$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'300.53"
$ws.Range("E2").Value = "'0.58%"
$ws.Range("D3").Value = "'31.59"
$ws.Range("E3").Value = "'0.81%"
$ws.Range("D4").Value = "'5.100"
$ws.Range("E4").Value = "'-0.78%"
$ws.Range("D5").Value = "'0.07815"
$ws.Range("E5").Value = "'-1.57%"
$ws.Range("D6").Value = "'2.332"
$ws.Range("E6").Value = "'-8.27%"
$ws.Range("D7").Value = "'7.801"
$ws.Range("E7").Value = "'-0.45%"
$ws.Range("D8").Value = "'3.832"
$ws.Range("E8").Value = "'0.16%"
$ws.Range("D9").Value = "'0.9153"
$ws.Range("E9").Value = "'0.91%"
$ws.Range("D10").Value = "'0.1761"
$ws.Range("E10").Value = "'1.41%"
$ws.Range("D11").Value = "'0.07533"
$ws.Range("E11").Value = "'4.01%"
$ws.Range("D12").Value = "'0.09137"
$ws.Range("E12").Value = "'13.48%"
$ws.Range("D13").Value = "'0.03094"
$ws.Range("E13").Value = "'2.38%"
$ws.Range("D15").Value = "'0.001511"
$ws.Range("E15").Value = "'1.11%"
$ws.Range("D16").Value = "'0.005798"
$ws.Range("E16").Value = "'-3.73%"
$ws.Range("D17").Value = "'3.482"
$ws.Range("E17").Value = "'-0.59%"
$ws.Range("E18").Value = "'-0.27%"
$ws.Range("E20").Value = "'1.53%"
$ws.Range("D21").Value = "'4.031"
$ws.Range("E21").Value = "'-12.87%"
$ws.Range("D23").Value = "'0.04599"
$ws.Range("E23").Value = "'0.37%"
$ws.Range("E24").Value = "'-0.56%"
$ws.Range("E25").Value = "'0.10%"
$ws.Range("D26").Value = "'0.0001250"
$ws.Range("E26").Value = "'6.02%"
$ws.Range("E27").Value = "'-1.35%"
$ws.Range("D39").Value = "'0.01775"
$ws.Range("E39").Value = "'-3.85%"
$ws.Range("D40").Value = "'0.04791"
$ws.Range("E40").Value = "'5.99%"
$ws.Range("D41").Value = "'0.007380"
$ws.Range("E41").Value = "'4.92%"
$ws.Range("D42").Value = "'0.1356"
$ws.Range("E42").Value = "'0.90%"
$ws.Range("D43").Value = "'0.002191"
$ws.Range("E43").Value = "'-2.15%"
$ws.Range("E44").Value = "'-1.66%"
$ws.Range("D45").Value = "'0.00006212"
$ws.Range("E45").Value = "'-3.67%"
$ws.Range("D46").Value = "'0.00000000750"
$ws.Range("E46").Value = "'0.09%"
$ws.Range("E47").Value = "'28.89%"
$ws.Range("D48").Value = "'0.7427"
$ws.Range("E48").Value = "'-9.49%"
$ws.Range("D49").Value = "'0.00002101"
$ws.Range("E49").Value = "'0.09%"
$ws.Range("D50").Value = "'0.0002001"
$ws.Range("E50").Value = "'0.09%"
